# Add a new worksheet "Akagi et al." (with a small cell-type/measurement
# table) to the workbook, after the existing "Chaudhary et al." sheet.

$wb = $excel.ActiveWorkbook

# Sheet that already carries the bold/bordered/centered style used for
# header cells and the numeric index column on every other sheet in this
# workbook - we reuse it below instead of re-creating the formatting from
# scratch so we don't end up with duplicate style records.
$styleSource = $wb.Worksheets.Item(1)

# Add the new worksheet right after the last existing sheet.
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Akagi et al."

# Header row.
$newSheet.Range("B1").Value = "cell type"
$newSheet.Range("C1").Value = "measurements"

# Data rows - index column, cell type name, number of measurements.
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "jurkat"
$newSheet.Cells.Item(2, 3).Value = 50

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "MEF"
$newSheet.Cells.Item(3, 3).Value = 50

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "hMSC"
$newSheet.Cells.Item(4, 3).Value = 50

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "hiPSC"
$newSheet.Cells.Item(5, 3).Value = 50

# Apply the same bold/border/center formatting used for the header row and
# the numeric index column on the other sheets, by copying the formats from
# an existing styled range (avoids creating duplicate style/font entries).
$styleSource.Range("B1:C1").Copy()
$newSheet.Range("B1:C1").PasteSpecial(-4122)

$styleSource.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$excel.CutCopyMode = 0
